$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (serial 45308 -> 45309)
$ws.Range("A1").Value = 45309

# Update unit prices in column D for rows 28-34
$ws.Range("D28").Value = 636.951
$ws.Range("D29").Value = 667.362
$ws.Range("D30").Value = 1075.196
$ws.Range("D31").Value = 1526.765
$ws.Range("D32").Value = 2133.272
$ws.Range("D33").Value = 2737.882
$ws.Range("D34").Value = 4277.943
